$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Data for each (sending, target) pair, in the row order needed (rows 2-10)
# Columns: E F G H I J K L M N O P Q R S T
$data = @(
    # Row 2: ECs -> ECs
    @(3, 1, 1.825549, 5.476647, 0.04696949406168958, 0.04696949406168958, 2, 0.6666666666666666, 19.86261233333333, 59.587837, 0.1710751304955294, 0.1710751304955294, 36.26017208250433, 326.341548742539, 0.008035312325912541, 0.008035312325912541),
    # Row 3: ECs -> FAPs
    @(3, 1, 1.825549, 5.476647, 0.04696949406168958, 0.04696949406168958, 3, 1, 69.67747766666666, 209.032433, 0.6001266794307873, 0.6001266794307873, 127.1996496769057, 1144.796847092151, 0.02818764650578585, 0.02818764650578585),
    # Row 4: ECs -> sCs
    @(3, 1, 1.825549, 5.476647, 0.04696949406168958, 0.04696949406168958, 3, 1, 26.564526, 79.693578, 0.2287981900736832, 0.2287981900736832, 48.494843874774, 436.453594872966, 0.01074653522999119, 0.01074653522999119),
    # Row 5: FAPs -> ECs
    @(3, 1, 23.57737633333333, 70.73212899999999, 0.6066215903701957, 0.6066215903701957, 2, 0.6666666666666666, 19.86261233333333, 59.587837, 0.1710751304955294, 0.1710751304955294, 468.3082859461081, 4214.774573514972, 0.1037778677339868, 0.1037778677339868),
    # Row 6: FAPs -> FAPs
    @(3, 1, 23.57737633333333, 70.73212899999999, 0.6066215903701957, 0.6066215903701957, 3, 1, 69.67747766666666, 209.032433, 0.6001266794307873, 0.6001266794307873, 1642.812112904428, 14785.30901613985, 0.3640498006998888, 0.3640498006998888),
    # Row 7: FAPs -> sCs
    @(3, 1, 23.57737633333333, 70.73212899999999, 0.6066215903701957, 0.6066215903701957, 3, 1, 26.564526, 79.693578, 0.2287981900736832, 0.2287981900736832, 626.3218266186179, 5636.896439567561, 0.13879392193632, 0.13879392193632),
    # Row 8: sCs -> ECs
    @(3, 1, 13.46376966666667, 40.391309, 0.3464089155681148, 0.3464089155681148, 2, 0.6666666666666666, 19.86261233333333, 59.587837, 0.1710751304955294, 0.1710751304955294, 267.4256374342926, 2406.830736908633, 0.05926195043563007, 0.05926195043563007),
    # Row 9: sCs -> FAPs
    @(3, 1, 13.46376966666667, 40.391309, 0.3464089155681148, 0.3464089155681148, 3, 1, 69.67747766666666, 209.032433, 0.6001266794307873, 0.6001266794307873, 938.1215102583106, 8443.093592324796, 0.2078892322251127, 0.2078892322251127),
    # Row 10: sCs -> sCs
    @(3, 1, 13.46376966666667, 40.391309, 0.3464089155681148, 0.3464089155681148, 3, 1, 26.564526, 79.693578, 0.2287981900736832, 0.2287981900736832, 357.658659368178, 3218.927934313602, 0.07925773290737199, 0.07925773290737201)
)

$rowPairs = @(
    @("ECs", "ECs"),
    @("ECs", "FAPs"),
    @("ECs", "sCs"),
    @("FAPs", "ECs"),
    @("FAPs", "FAPs"),
    @("FAPs", "sCs"),
    @("sCs", "ECs"),
    @("sCs", "FAPs"),
    @("sCs", "sCs")
)

for ($i = 0; $i -lt $rowPairs.Count; $i++) {
    $r = $i + 2
    $sending = $rowPairs[$i][0]
    $target = $rowPairs[$i][1]

    $ws.Cells.Item($r, 1).Value = $sending
    $ws.Cells.Item($r, 2).Value = "Tgfb3"
    $ws.Cells.Item($r, 3).Value = "Tgfbr3"
    $ws.Cells.Item($r, 4).Value = $target

    $vals = $data[$i]
    for ($c = 0; $c -lt $vals.Count; $c++) {
        $ws.Cells.Item($r, $c + 5).Value = $vals[$c]
    }
}
